# Auto-generated edit script: apply F-column (想去人数) value updates
# per the commit diff, matched by (sheet, row, expected-old-value) for safety.
$wb = $excel.ActiveWorkbook

$updates = @{}
$updates["展览"] = @(
    @{ Row = 6; Old = 168; New = 170 },
    @{ Row = 7; Old = 163; New = 165 },
    @{ Row = 8; Old = 4175; New = 4182 },
    @{ Row = 10; Old = 57; New = 58 },
    @{ Row = 11; Old = 172; New = 173 },
    @{ Row = 13; Old = 6005; New = 6021 },
    @{ Row = 14; Old = 787; New = 788 },
    @{ Row = 16; Old = 2317; New = 2321 },
    @{ Row = 19; Old = 463; New = 464 },
    @{ Row = 20; Old = 9039; New = 9071 },
    @{ Row = 21; Old = 39; New = 40 },
    @{ Row = 22; Old = 2436; New = 2442 },
    @{ Row = 24; Old = 2297; New = 2299 },
    @{ Row = 25; Old = 2416; New = 2422 },
    @{ Row = 26; Old = 1383; New = 1384 },
    @{ Row = 28; Old = 1945; New = 1948 },
    @{ Row = 34; Old = 277; New = 278 },
    @{ Row = 36; Old = 42; New = 46 },
    @{ Row = 37; Old = 32; New = 33 },
    @{ Row = 41; Old = 93; New = 94 },
    @{ Row = 42; Old = 230; New = 232 },
    @{ Row = 43; Old = 1512; New = 1515 },
    @{ Row = 44; Old = 2471; New = 2477 },
    @{ Row = 45; Old = 912; New = 914 },
    @{ Row = 47; Old = 1248; New = 1249 },
    @{ Row = 48; Old = 13; New = 14 }
)

$updates["演出"] = @(
    @{ Row = 22; Old = 48; New = 60 },
    @{ Row = 23; Old = 48; New = 60 }
)

$updates["本地生活"] = @(
    @{ Row = 2; Old = 689; New = 691 }
)

$updates["全部类型"] = @(
    @{ Row = 4; Old = 689; New = 691 },
    @{ Row = 9; Old = 168; New = 170 },
    @{ Row = 12; Old = 163; New = 165 },
    @{ Row = 13; Old = 4175; New = 4182 },
    @{ Row = 14; Old = 57; New = 58 },
    @{ Row = 15; Old = 172; New = 173 },
    @{ Row = 16; Old = 6005; New = 6021 },
    @{ Row = 17; Old = 787; New = 788 },
    @{ Row = 19; Old = 2317; New = 2321 },
    @{ Row = 21; Old = 463; New = 464 },
    @{ Row = 22; Old = 9039; New = 9071 },
    @{ Row = 24; Old = 2436; New = 2442 },
    @{ Row = 25; Old = 2297; New = 2299 },
    @{ Row = 26; Old = 1383; New = 1384 },
    @{ Row = 28; Old = 1945; New = 1948 },
    @{ Row = 33; Old = 277; New = 278 },
    @{ Row = 35; Old = 42; New = 46 },
    @{ Row = 36; Old = 32; New = 33 },
    @{ Row = 39; Old = 93; New = 94 },
    @{ Row = 40; Old = 230; New = 232 },
    @{ Row = 41; Old = 1512; New = 1515 },
    @{ Row = 42; Old = 2471; New = 2477 },
    @{ Row = 43; Old = 912; New = 914 },
    @{ Row = 48; Old = 1248; New = 1249 },
    @{ Row = 49; Old = 13; New = 14 },
    @{ Row = 50; Old = 48; New = 60 }
)

$mismatches = 0
foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $updates[$sheetName]) {
        $cell = $ws.Cells.Item($u.Row, 6)
        $current = $cell.Value()
        if ($current -ne $u.Old) {
            $mismatches = $mismatches + 1
        }
        $cell.Value = $u.New
    }
}
Write-Host "mismatches:" $mismatches